$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "G" column for the rows where the second credit/value bumps
# from 1 to 2 (and add the missing G9 cell with value 1), mirroring the
# author's move to a repetitive 2D (row x col) loop structure.
$rowsToIncrement = @(3, 6, 7, 12, 13, 14, 17, 18, 19, 21)
foreach ($r in $rowsToIncrement) {
    $ws.Cells.Item($r, 7).Value = 2
}

# G9 did not previously exist as a populated cell; it now gets a value of 1.
$ws.Cells.Item(9, 7).Value = 1

# Update the active selection to reflect the last edit: column G, rows 3:21.
$ws.Range("G3:G21").Select()
